# Clear out the personal/sample data that had been filled into the
# "Imposto de Renda" template, restoring it back to a blank template.

$wb = $excel.ActiveWorkbook

$wsTitular  = $wb.Worksheets.Item("TITULAR")
$wsInformes = $wb.Worksheets.Item("INFORMES")
$wsNotas    = $wb.Worksheets.Item("NOTAS")

# --- TITULAR: clear every answer in column D (rows 6-19) and drop the
#     e-mail hyperlink that was attached to D16.
$wsTitular.Range("D6:D19").ClearContents()
$wsTitular.Hyperlinks.Delete()

# --- INFORMES: clear the bank / value / attachment fields for all three
#     "informes" blocks.
$wsInformes.Range("D10:D12").ClearContents()
$wsInformes.Range("D15:D17").ClearContents()
$wsInformes.Range("D20:D22").ClearContents()

# --- NOTAS: clear the single sample entry in the table.
$wsNotas.Range("C9:E9").ClearContents()

# --- Restore view state to the blank-template defaults: TITULAR active,
#     each sheet's selection back on its first input cell.
$wsInformes.Select()
$wsInformes.Range("D10").Select()

$wsNotas.Select()
$wsNotas.Range("C9").Select()

$wsTitular.Select()
$wsTitular.Range("D6").Select()

Write-Output "edit applied"
